$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.170.70'
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").Value = '1.752.25'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.36'
$ws.Range("E5").Value = '  +1.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5289'
$ws.Range("E7").Value = '  +2.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2827'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06198'
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("D10").Value = '1.744.31'
$ws.Range("E10").Value = '  -0.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07213'
$ws.Range("E11").Value = '  +3.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.55'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6488'
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.641'
$ws.Range("E14").Value = '  +2.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '78.74'
$ws.Range("E15").Value = '  +2.12%  '
$ws.Range("E16").Value = '  +0.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9996'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").Value = '26.062.39'
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("E19").Value = '  +2.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006758'
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("D21").Value = '1.972.10'
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.343'
$ws.Range("E22").Value = '  +4.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.781'
$ws.Range("E23").Value = '  +2.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.252'
$ws.Range("E24").Value = '  +1.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.08'
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.517'
$ws.Range("E26").Value = '  +0.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.37'
$ws.Range("E27").Value = '  +1.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.818'
$ws.Range("E28").Value = '  -0.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '105.66'
$ws.Range("E29").Value = '  +2.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08328'
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.815'
$ws.Range("E31").Value = '  +4.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.663'
$ws.Range("E32").Value = '  +6.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04623'
$ws.Range("E33").Value = '  +4.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.650'
$ws.Range("E34").Value = '  +1.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.021'
$ws.Range("E35").Value = '  +3.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6372'
$ws.Range("E36").Value = '  +4.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.704'
$ws.Range("E37").Value = '  +0.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01626'
$ws.Range("E38").Value = '  +3.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.991'
$ws.Range("E39").Value = '  +2.21%  '
$ws.Range("E40").Value = '  +0.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.74'
$ws.Range("E41").Value = '  +1.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.3951'
$ws.Range("E42").Value = '  +1.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7513'
$ws.Range("E43").Value = '  +2.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.071'
$ws.Range("E44").Value = '  +2.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1157'
$ws.Range("E45").Value = '  +3.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.396'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05351'
$ws.Range("E47").Value = '  -1.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '31.16'
$ws.Range("E48").Value = '  +4.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.67'
$ws.Range("E49").Value = '  +3.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3494'
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.622'
$ws.Range("E51").Value = '  +1.50%  '
